# Booking sheet update: fix checkout dates / values (update, delete, get
# booking endpoints).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "checkout " -> "checkout" (trailing space trimmed)
$ws.Range("F1").Value = "checkout"

# Row 2: checkout text corrected to match checkin (2018-01-01)
$ws.Range("F2").Value = "2018-01-01 "

# Row 3 & 4: checkout columns become real dates (instead of placeholder
# text). Copy the existing date formatting from E5 (numFmtId 14) so the
# cells pick up the same style as the other date cells in the sheet.
$ws.Range("E5").Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4122) | Out-Null
$ws.Range("F3").Value = 43112

$ws.Range("E5").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null
$ws.Range("F4").Value = 43110

# Row 5: correct both checkin and checkout dates
$ws.Range("E5").Value = 43110
$ws.Range("F5").Value = 43111

# Row 6: typo fix in checkout text
$ws.Range("F6").Value = "2019-14-02"

# Restore selection as recorded in the workbook
$ws.Range("F12").Select() | Out-Null
